$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 17 by copying row 16 (same worker, new period) - this
# shifts rows 17-22 down to 18-23, duplicating formatting.
$ws.Rows.Item(16).Copy()
$ws.Rows.Item(17).Insert()

# New row: same worker/doc/name, but period 2509 instead of 2508
$ws.Range("E17").Value = "2509"

# Update totals: two periods in mora now (2508 + 2509), so the mora
# value doubles and the period count increases to 2.
$ws.Range("E11").Value = 113880
$ws.Range("F13").Value = 2
